# Generate Report for Handback
# Updates the "Latest HO Xliff Generate Date" on the Overview sheet, and
# the "Correspond Handoff Datetime" / "Correspond Handback DateTime"
# columns on the zh-cn / de-de report sheets for the second data row
# (the 89b0dd29-... file), reflecting a newer handback run.

$wb = $excel.ActiveWorkbook

# --- Overview sheet --------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
# Row 3 corresponds to 89b0dd29-cb79-48f1-9fdc-c0e62a8961ea.md
$overview.Range("G3").Value = "2016-10-18 11:45:31"

# --- zh-cn sheet -------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
# Row 3 corresponds to 89b0dd29-cb79-48f1-9fdc-c0e62a8961ea.md
$zhcn.Range("H3").Value = "2016-10-18 11:45:20"
$zhcn.Range("K3").Value = "2016-10-18 11:46:02"

# --- de-de sheet -------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
# Row 3 corresponds to 89b0dd29-cb79-48f1-9fdc-c0e62a8961ea.md
$dede.Range("H3").Value = "2016-10-18 11:45:31"
$dede.Range("K3").Value = "2016-10-18 11:46:19"
